$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "H2"  = 84
    "I2"  = 242
    "J2"  = 1082
    "K2"  = 2
    "L2"  = 311
    "M2"  = 24
    "N2"  = 162
    "O2"  = 0
    "P2"  = 6
    "Q2"  = 6
    "R2"  = 14
    "S2"  = 110
    "T2"  = 188
    "U2"  = 16
    "V2"  = 1655
    "W2"  = 1
    "X2"  = 1656
    "Y2"  = 1
    "Z2"  = 30
    "AA2" = 14
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
